$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36 (Idaho) now failed with a timeout instead of succeeding.
# Clear the previously-populated numeric columns B:H back to blank/empty
# string cells (also drops the date style that was on B36), flip the
# "Pct Includes Hispanic Black" boolean to False, and replace the status
# message with the timeout error text.

$cols = @("B", "C", "D", "E", "F", "G", "H")
foreach ($col in $cols) {
    $cell = $ws.Range($col + "36")
    $cell.Value = "'"
    $cell.ClearFormats()
}

$ws.Range("J36").Value = $false

$ws.Range("O36").Value = "An error occurred. ... TimeoutException('', None, None)"
